# The stimuli for the congruent-practice list were converted from .bmp to
# .png files, so the stim_path column (D2:D9) needs to reflect the new
# file extension.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "stimuli/1.png"
$ws.Range("D3").Value = "stimuli/5.png"
$ws.Range("D4").Value = "stimuli/7.png"
$ws.Range("D5").Value = "stimuli/9.png"
$ws.Range("D6").Value = "stimuli/14.png"
$ws.Range("D7").Value = "stimuli/15.png"
$ws.Range("D8").Value = "stimuli/19.png"
$ws.Range("D9").Value = "stimuli/20.png"

# Update the sheet's last active/selected cell to D2 (as saved in the
# workbook's sheetView), matching what was recorded after the edit.
$ws.Range("D2").Select()
